$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 96
$ws.Cells.Item(96, 8).Value = 875.75
$ws.Cells.Item(96, 9).Value = 572
$ws.Cells.Item(96, 10).Value = 1058
$ws.Cells.Item(96, 11).Value = 1716
$ws.Cells.Item(96, 12).Value = 3174
$ws.Cells.Item(96, 13).Value = -343
$ws.Cells.Item(96, 14).Value = -5920
# Row 111
$ws.Cells.Item(111, 8).Value = 517.375
$ws.Cells.Item(111, 9).Value = 448
$ws.Cells.Item(111, 10).Value = 670
$ws.Cells.Item(111, 11).Value = 1344
$ws.Cells.Item(111, 12).Value = 2010
$ws.Cells.Item(111, 13).Value = 1723
$ws.Cells.Item(111, 14).Value = -8144
# Row 129
$ws.Cells.Item(129, 8).Value = 983.2033699999999
$ws.Cells.Item(129, 9).Value = 335
$ws.Cells.Item(129, 10).Value = 1030.3455
$ws.Cells.Item(129, 11).Value = 1005
$ws.Cells.Item(129, 12).Value = 3091.0365
$ws.Cells.Item(129, 13).Value = 3995
$ws.Cells.Item(129, 14).Value = -13091.0365
# Row 131
$ws.Cells.Item(131, 8).Value = 4738.32
$ws.Cells.Item(131, 9).Value = 4903.391
$ws.Cells.Item(131, 10).Value = 2840
$ws.Cells.Item(131, 11).Value = 14710.173
$ws.Cells.Item(131, 12).Value = 8520
$ws.Cells.Item(131, 13).Value = -9670.172999999999
$ws.Cells.Item(131, 14).Value = -18600
# Row 138
$ws.Cells.Item(138, 8).Value = 2484.7795
$ws.Cells.Item(138, 9).Value = 963.5
$ws.Cells.Item(138, 10).Value = 4058.5173
$ws.Cells.Item(138, 11).Value = 2890.5
$ws.Cells.Item(138, 12).Value = 12175.5519
$ws.Cells.Item(138, 13).Value = 2249.5
$ws.Cells.Item(138, 14).Value = -22455.5519

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 1442.7
$ws.Cells.Item(2, 9).Value = 845.25
$ws.Cells.Item(2, 10).Value = 3832.5
$ws.Cells.Item(2, 11).Value = 845.25
$ws.Cells.Item(2, 12).Value = 3832.5
$ws.Cells.Item(2, 13).Value = -732.25
$ws.Cells.Item(2, 14).Value = -4058.5
# Row 97
$ws.Cells.Item(97, 8).Value = 977.4838999999999
$ws.Cells.Item(97, 9).Value = 764.6667
$ws.Cells.Item(97, 10).Value = 1424.4
$ws.Cells.Item(97, 11).Value = 764.6667
$ws.Cells.Item(97, 12).Value = 1424.4
$ws.Cells.Item(97, 13).Value = -268.6667
$ws.Cells.Item(97, 14).Value = -2416.4
# Row 102
$ws.Cells.Item(102, 8).Value = 1992.2
$ws.Cells.Item(102, 9).Value = 1855.3549
$ws.Cells.Item(102, 10).Value = 3052.75
$ws.Cells.Item(102, 11).Value = 1855.3549
$ws.Cells.Item(102, 12).Value = 3052.75
$ws.Cells.Item(102, 13).Value = -233.3549
$ws.Cells.Item(102, 14).Value = -6296.75
# Row 116
$ws.Cells.Item(116, 8).Value = 1442.7
$ws.Cells.Item(116, 9).Value = 845.25
$ws.Cells.Item(116, 10).Value = 3832.5
$ws.Cells.Item(116, 11).Value = 845.25
$ws.Cells.Item(116, 12).Value = 3832.5
$ws.Cells.Item(116, 13).Value = 1448.75
$ws.Cells.Item(116, 14).Value = -8420.5
# Row 122
$ws.Cells.Item(122, 8).Value = 903.75
$ws.Cells.Item(122, 9).Value = 907.93335
$ws.Cells.Item(122, 10).Value = 891.2
$ws.Cells.Item(122, 11).Value = 2723.80005
$ws.Cells.Item(122, 12).Value = 2673.6
$ws.Cells.Item(122, 13).Value = -273.8000499999998
$ws.Cells.Item(122, 14).Value = -7573.6
# Row 132
$ws.Cells.Item(132, 8).Value = 19525.7
$ws.Cells.Item(132, 9).Value = 20391.148
$ws.Cells.Item(132, 10).Value = 14333
$ws.Cells.Item(132, 11).Value = 61173.444
$ws.Cells.Item(132, 12).Value = 42999
$ws.Cells.Item(132, 13).Value = -58643.444
$ws.Cells.Item(132, 14).Value = -48059

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 1442.7
$ws.Cells.Item(3, 9).Value = 845.25
$ws.Cells.Item(3, 10).Value = 3832.5
$ws.Cells.Item(3, 11).Value = 845.25
$ws.Cells.Item(3, 12).Value = 3832.5
$ws.Cells.Item(3, 13).Value = -731.25
$ws.Cells.Item(3, 14).Value = -4060.5
# Row 94
$ws.Cells.Item(94, 8).Value = 1415.381
$ws.Cells.Item(94, 9).Value = 982.25
$ws.Cells.Item(94, 10).Value = 1992.8889
$ws.Cells.Item(94, 11).Value = 982.25
$ws.Cells.Item(94, 12).Value = 1992.8889
$ws.Cells.Item(94, 13).Value = -531.25
$ws.Cells.Item(94, 14).Value = -2894.8889
# Row 134
$ws.Cells.Item(134, 8).Value = 26132.146
$ws.Cells.Item(134, 9).Value = 33646.03
$ws.Cells.Item(134, 10).Value = 2839.1
$ws.Cells.Item(134, 11).Value = 100938.09
$ws.Cells.Item(134, 12).Value = 8517.299999999999
$ws.Cells.Item(134, 13).Value = -98403.09
$ws.Cells.Item(134, 14).Value = -13587.3

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 40716.59
$ws.Cells.Item(31, 9).Value = 60585.363
$ws.Cells.Item(31, 10).Value = 4290.5
$ws.Cells.Item(31, 11).Value = 60585.363
$ws.Cells.Item(31, 12).Value = 4290.5
$ws.Cells.Item(31, 13).Value = -60290.363
$ws.Cells.Item(31, 14).Value = -4880.5
# Row 34
$ws.Cells.Item(34, 8).Value = 40716.59
$ws.Cells.Item(34, 9).Value = 60585.363
$ws.Cells.Item(34, 10).Value = 4290.5
$ws.Cells.Item(34, 11).Value = 60585.363
$ws.Cells.Item(34, 12).Value = 4290.5
$ws.Cells.Item(34, 13).Value = -60383.363
$ws.Cells.Item(34, 14).Value = -4694.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 2231.818
$ws.Cells.Item(5, 10).Value = 2318.3333
$ws.Cells.Item(5, 12).Value = 6954.999899999999
$ws.Cells.Item(5, 14).Value = -7178.999899999999
# Row 9
$ws.Cells.Item(9, 8).Value = 1440.909
$ws.Cells.Item(9, 10).Value = 1575
$ws.Cells.Item(9, 12).Value = 4725
$ws.Cells.Item(9, 14).Value = -5173
# Row 75
$ws.Cells.Item(75, 8).Value = 1763.6666
$ws.Cells.Item(75, 10).Value = 1763.6666
$ws.Cells.Item(75, 12).Value = 5290.9998
$ws.Cells.Item(75, 14).Value = -7286.9998
# Row 78
$ws.Cells.Item(78, 8).Value = 1763.6666
$ws.Cells.Item(78, 10).Value = 1763.6666
$ws.Cells.Item(78, 12).Value = 15872.9994
$ws.Cells.Item(78, 14).Value = -25856.9994
# Row 122
$ws.Cells.Item(122, 8).Value = 1700
$ws.Cells.Item(122, 9).Value = 657.1429000000001
$ws.Cells.Item(122, 10).Value = 2363.6365
$ws.Cells.Item(122, 11).Value = 5914.2861
$ws.Cells.Item(122, 12).Value = 21272.7285
$ws.Cells.Item(122, 13).Value = -3464.2861
$ws.Cells.Item(122, 14).Value = -26172.7285
# Row 135
$ws.Cells.Item(135, 8).Value = 2231.818
$ws.Cells.Item(135, 10).Value = 2318.3333
$ws.Cells.Item(135, 12).Value = 20864.9997
$ws.Cells.Item(135, 14).Value = -25934.9997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Cells.Item(113, 8).Value = 1100
$ws.Cells.Item(113, 9).Value = 970
$ws.Cells.Item(113, 10).Value = 1178
$ws.Cells.Item(113, 11).Value = 970
$ws.Cells.Item(113, 12).Value = 1178
$ws.Cells.Item(113, 13).Value = 1200
$ws.Cells.Item(113, 14).Value = -5518

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Cells.Item(82, 8).Value = 2212.475
$ws.Cells.Item(82, 9).Value = 1433.25
$ws.Cells.Item(82, 10).Value = 2546.4285
$ws.Cells.Item(82, 11).Value = 1433.25
$ws.Cells.Item(82, 12).Value = 2546.4285
$ws.Cells.Item(82, 13).Value = -1072.25
$ws.Cells.Item(82, 14).Value = -3268.4285
# Row 85
$ws.Cells.Item(85, 8).Value = 2212.475
$ws.Cells.Item(85, 9).Value = 1433.25
$ws.Cells.Item(85, 10).Value = 2546.4285
$ws.Cells.Item(85, 11).Value = 1433.25
$ws.Cells.Item(85, 12).Value = 2546.4285
$ws.Cells.Item(85, 13).Value = -185.25
$ws.Cells.Item(85, 14).Value = -5042.4285
# Row 88
$ws.Cells.Item(88, 8).Value = 25099.5
$ws.Cells.Item(88, 10).Value = 25099.5
$ws.Cells.Item(88, 12).Value = 25099.5
$ws.Cells.Item(88, 14).Value = -25955.5
# Row 91
$ws.Cells.Item(91, 8).Value = 25099.5
$ws.Cells.Item(91, 10).Value = 25099.5
$ws.Cells.Item(91, 12).Value = 25099.5
$ws.Cells.Item(91, 14).Value = -28063.5
# Row 122
$ws.Cells.Item(122, 8).Value = 1828
$ws.Cells.Item(122, 9).Value = 1848.909
$ws.Cells.Item(122, 10).Value = 1782
$ws.Cells.Item(122, 11).Value = 5546.727000000001
$ws.Cells.Item(122, 12).Value = 5346
$ws.Cells.Item(122, 13).Value = -3096.727000000001
$ws.Cells.Item(122, 14).Value = -10246

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 3800
$ws.Cells.Item(62, 10).Value = 3800
$ws.Cells.Item(62, 12).Value = 3800
$ws.Cells.Item(62, 14).Value = -5048
# Row 65
$ws.Cells.Item(65, 8).Value = 3800
$ws.Cells.Item(65, 10).Value = 3800
$ws.Cells.Item(65, 12).Value = 19000
$ws.Cells.Item(65, 14).Value = -25240
# Row 86
$ws.Cells.Item(86, 8).Value = 7193.75
$ws.Cells.Item(86, 10).Value = 7193.75
$ws.Cells.Item(86, 12).Value = 7193.75
$ws.Cells.Item(86, 14).Value = -9439.75
# Row 88
$ws.Cells.Item(88, 8).Value = 20000
$ws.Cells.Item(88, 9).Value = 20000
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 20000
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = -19594
$ws.Cells.Item(88, 14).ClearContents()
# Row 89
$ws.Cells.Item(89, 8).Value = 7193.75
$ws.Cells.Item(89, 10).Value = 7193.75
$ws.Cells.Item(89, 12).Value = 35968.75
$ws.Cells.Item(89, 14).Value = -47200.75
# Row 91
$ws.Cells.Item(91, 8).Value = 20000
$ws.Cells.Item(91, 9).Value = 20000
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 20000
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 13).Value = -18596
$ws.Cells.Item(91, 14).ClearContents()
# Row 92
$ws.Cells.Item(92, 8).Value = 39033.332
$ws.Cells.Item(92, 10).Value = 39033.332
$ws.Cells.Item(92, 12).Value = 39033.332
$ws.Cells.Item(92, 14).Value = -44025.332
# Row 96
$ws.Cells.Item(96, 8).Value = 3000000
$ws.Cells.Item(96, 9).Value = 1000000
$ws.Cells.Item(96, 10).Value = 5000000
$ws.Cells.Item(96, 11).Value = 1000000
$ws.Cells.Item(96, 12).Value = 5000000
$ws.Cells.Item(96, 13).Value = -998627
$ws.Cells.Item(96, 14).Value = -5002746
# Row 104
$ws.Cells.Item(104, 8).Value = 19273.6
$ws.Cells.Item(104, 10).Value = 19273.6
$ws.Cells.Item(104, 12).Value = 19273.6
$ws.Cells.Item(104, 14).Value = -26261.6
# Row 105
$ws.Cells.Item(105, 8).Value = 42711.25
$ws.Cells.Item(105, 10).Value = 42711.25
$ws.Cells.Item(105, 12).Value = 42711.25
$ws.Cells.Item(105, 14).Value = -49699.25

Write-Output "Yojimbo_Profits sheets updated successfully."